$d = $word.ActiveDocument

# Locate the bibliography entry that must be kept ("...Thomson Pioneira
# (2008).") and the trailing "Creative Commons Attribution" paragraph that
# marks the end of the footer block to be removed. Search by text instead
# of hard-coded paragraph indices so the script is resilient to minor
# structural differences.
$count = $d.Paragraphs.Count
$anchorIdx = -1
$tailIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Thomson Pioneira*") { $anchorIdx = $i }
    if ($t -like "*Creative Commons Attribution*") { $tailIdx = $i }
}

# Remove everything from just after the bibliography entry's paragraph mark
# through the end of the copyright/footer paragraph (inclusive of its
# paragraph mark). That deletes: the blank paragraph, the "Ver no Jupiter
# Salvar em pdf Salvar em docx" paragraph, and the "© 2020 ... Creative
# Commons Attribution" paragraph, while leaving the following blank
# paragraph and the page-break paragraph untouched.
$startPos = $d.Paragraphs.Item($anchorIdx).Range.End
$endPos = $d.Paragraphs.Item($tailIdx).Range.End

$toDelete = $d.Range($startPos, $endPos)
$toDelete.Delete()
